$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.089.56"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.426.58"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "2.413.25"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "2.871.54"
$ws.Range("D17").Value = "60.862.39"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "2.414.08"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "590.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.02%  "
$ws.Range("D29").Value = "2.541.21"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "0.0₃0939"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.80%  "
$ws.Range("E45").Value = "  -4.44%  "
$ws.Range("D46").Value = "0.0₆0292"
$ws.Range("E46").Value = "  +14.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0505"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.28%  "
